$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 3).Value = "OUT -> 2017/02/10 07:42"
$ws.Cells.Item(4, 4).Value = "OUT -> 2017/02/14 14:11"
$ws.Cells.Item(5, 1).Value = "IN -> 2017/02/14 08:23"
$ws.Cells.Item(5, 3).Value = "IN -> 2017/02/10 18:03"
$ws.Cells.Item(5, 4).Value = "IN -> 2017/02/16 14:29"
$ws.Cells.Item(5, 5).Value = "IN -> 2017/02/15 13:02"
$ws.Cells.Item(6, 1).Value = "OUT -> 2017/02/15 13:37"
$ws.Cells.Item(6, 3).Value = "OUT -> 2017/02/14 16:29"
$ws.Cells.Item(6, 4).Value = "OUT -> 2017/02/23 11:46"
$ws.Cells.Item(6, 5).Value = "OUT -> 2017/02/16 16:39"
$ws.Cells.Item(7, 1).Value = "IN -> 2017/02/16 10:43"
$ws.Cells.Item(7, 3).Value = "IN -> 2017/02/15 14:03"
$ws.Cells.Item(7, 5).Value = "IN -> 2017/02/22 20:10"
$ws.Cells.Item(8, 1).Value = "OUT -> 2017/02/16 14:29"
$ws.Cells.Item(8, 3).Value = "OUT -> 2017/02/15 20:50"
$ws.Cells.Item(8, 5).Value = "OUT -> 2017/02/23 17:26"
$ws.Cells.Item(9, 1).Value = "IN -> 2017/02/28 08:07"
$ws.Cells.Item(9, 3).Value = "IN -> 2017/02/17 07:53"
$ws.Cells.Item(9, 5).Value = "IN -> 2017/02/28 14:56"
$ws.Cells.Item(10, 1).Value = "OUT -> 2017/03/01 14:13"
$ws.Cells.Item(10, 3).Value = "OUT -> 2017/02/21 14:04"
$ws.Cells.Item(11, 3).Value = "IN -> 2017/02/21 16:26"
$ws.Cells.Item(12, 3).Value = "OUT -> 2017/02/21 19:01"
$ws.Cells.Item(13, 3).Value = "IN -> 2017/02/22 13:42"
$ws.Cells.Item(14, 3).Value = "OUT -> 2017/02/22 20:50"
$ws.Cells.Item(15, 3).Value = "IN -> 2017/02/22 20:50"
$ws.Cells.Item(16, 3).Value = "OUT -> 2017/02/22 20:50"
$ws.Cells.Item(17, 3).Value = "IN -> 2017/02/22 20:50"
$ws.Cells.Item(18, 3).Value = "OUT -> 2017/02/24 07:57"
$ws.Cells.Item(19, 3).Value = "IN -> 2017/02/24 17:54"
$ws.Cells.Item(20, 3).Value = "OUT -> 2017/02/24 17:54"
$ws.Cells.Item(21, 3).Value = "IN -> 2017/02/24 17:54"
$ws.Cells.Item(22, 3).Value = "OUT -> 2017/02/27 20:30"
$ws.Cells.Item(23, 3).Value = "IN -> 2017/02/28 14:16"
$ws.Cells.Item(24, 3).Value = "OUT -> 2017/02/28 16:26"
$ws.Cells.Item(25, 3).Value = "IN -> 2017/02/28 16:27"
$ws.Cells.Item(26, 3).Value = "OUT -> 2017/03/01 14:27"
$ws.Cells.Item(27, 3).Value = "IN -> 2017/03/01 18:06"
$ws.Cells.Item(28, 3).Value = "OUT -> 2017/03/01 20:47"
$ws.Cells.Item(29, 3).Value = "IN -> 2017/03/02 16:24"
$ws.Cells.Item(30, 3).Value = "OUT -> 2017/03/02 16:25"
$ws.Cells.Item(31, 3).Value = "IN -> 2017/03/03 11:13"
